{"js": "// Fix the meeting date: the month was mistakenly typed as \"abril\" (April)\n// and must read \"marzo\" (March) instead, e.g. \"Fecha: 12 de marzo de 2025\".\nconst body = context.document.body;\n\nconst results = body.search(\"Fecha: 12 de abril de 2025\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target date text not found\");\n}\n\nfor (const range of results.items) {\n  range.insertText(\"Fecha: 12 de marzo de 2025\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix the meeting date: the month was mistakenly typed as \"abril\" (April)\n# and must read \"marzo\" (March) instead, e.g. \"Fecha: 12 de marzo de 2025\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Fecha: 12 de abril de 2025\"\n$find.Replacement.Text = \"Fecha: 12 de marzo de 2025\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$found = $find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $find.MatchSoundsLike,\n    $find.MatchAllWordForms,\n    $find.Forward,\n    $find.Wrap,\n    $find.Format,\n    $find.Replacement.Text,\n    2\n)\n\nif (-not $found) {\n    throw \"Target date text not found\"\n}\n"}
